# Minor changes to Fantasy setup
#
# 1) "All Fun little Glory" roster: Brad Brock (8,4) is replaced by
#    Sean Cotnam (5,6) on row 4.
# 2) Summary sheet: because "All Fun little Glory"'s total drops from 99
#    to 98 points, it now sorts below "Suzanne's Quick Finishers" (still
#    99) in the points-descending leaderboard, so rows 6 and 7 swap.

$wb = $excel.ActiveWorkbook

# --- 1) Update the "All Fun little Glory" roster row ---
$team = $wb.Worksheets.Item("All Fun little Glory")
$team.Range("A4").Value = "Sean Cotnam"
$team.Range("B4").Value = 5
$team.Range("C4").Value = 6

# --- 2) Update the Summary leaderboard ordering ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("A6").Value = "Suzanne's Quick Finishers"
$summary.Range("B6").Value = 99
$summary.Range("A7").Value = "All Fun little Glory"
$summary.Range("B7").Value = 98
